$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated result values (row, columnIndex, new value)
$updates = @(
    @(4, 2, 0.3),
    @(4, 5, 0.18),
    @(4, 6, 0.028),
    @(4, 7, 0.167),
    @(4, 8, 0.214),
    @(4, 10, 0.157),
    @(4, 11, 0.321),
    @(4, 12, 0.101),
    @(4, 13, 0.318),
    @(4, 14, 0.256),
    @(4, 16, 0.148),
    @(4, 17, 0.512),
    @(4, 18, 0.22),
    @(4, 19, 0.469),
    @(4, 20, 0.262),
    @(4, 21, 0.08500000000000001),
    @(4, 23, 0.243),
    @(4, 24, 0.044),
    @(4, 25, 0.209),
    @(4, 26, 0.454),
    @(4, 27, 0.13),
    @(4, 28, 0.361),
    @(4, 31, 0.081),
    @(4, 32, 0.71),
    @(4, 33, 0.104),
    @(4, 34, 0.323),
    @(4, 35, 0.643),
    @(4, 36, 0.17),
    @(4, 37, 0.412),
    @(4, 38, 0.66),
    @(4, 40, 0.351),
    @(4, 41, 0.671),
    @(5, 2, 0.833),
    @(5, 3, 0.139),
    @(5, 4, 0.373),
    @(5, 5, 0.714),
    @(5, 6, 0.204),
    @(5, 7, 0.452),
    @(5, 8, 0.857),
    @(5, 9, 0.122),
    @(5, 10, 0.35),
    @(5, 11, 0.619),
    @(5, 12, 0.236),
    @(5, 13, 0.486),
    @(5, 14, 0.8100000000000001),
    @(5, 15, 0.154),
    @(5, 16, 0.393),
    @(5, 17, 0.571),
    @(5, 18, 0.245),
    @(5, 19, 0.495),
    @(5, 20, 0.548),
    @(5, 21, 0.248),
    @(5, 22, 0.498),
    @(5, 23, 0.738),
    @(5, 24, 0.193),
    @(5, 25, 0.44),
    @(5, 26, 0.833),
    @(5, 27, 0.139),
    @(5, 28, 0.373),
    @(5, 29, 0.738),
    @(5, 30, 0.193),
    @(5, 31, 0.44),
    @(5, 32, 0.952),
    @(5, 33, 0.045),
    @(5, 34, 0.213),
    @(5, 35, 0.762),
    @(5, 36, 0.181),
    @(5, 37, 0.426),
    @(5, 38, 0.905),
    @(5, 39, 0.08599999999999999),
    @(5, 40, 0.294),
    @(5, 41, 0.873),
    @(6, 2, 0.441),
    @(6, 5, 0.288),
    @(6, 8, 0.342),
    @(6, 11, 0.423),
    @(6, 14, 0.389),
    @(6, 17, 0.54),
    @(6, 20, 0.355),
    @(6, 23, 0.366),
    @(6, 26, 0.588),
    @(6, 32, 0.8129999999999999),
    @(6, 35, 0.697),
    @(6, 38, 0.763),
    @(6, 41, 0.758),
    @(7, 2, 0.615),
    @(7, 5, 0.448),
    @(7, 8, 0.535),
    @(7, 11, 0.522),
    @(7, 14, 0.5649999999999999),
    @(7, 17, 0.5580000000000001),
    @(7, 20, 0.45),
    @(7, 23, 0.524),
    @(7, 26, 0.714),
    @(7, 29, 0.367),
    @(7, 32, 0.891),
    @(7, 35, 0.735),
    @(7, 38, 0.842),
    @(7, 41, 0.823),
    @(8, 2, 0.768),
    @(8, 3, 0.144),
    @(8, 4, 0.379),
    @(8, 5, 0.603),
    @(8, 7, 0.426),
    @(8, 8, 0.747),
    @(8, 9, 0.132),
    @(8, 10, 0.364),
    @(8, 11, 0.545),
    @(8, 12, 0.208),
    @(8, 13, 0.456),
    @(8, 14, 0.728),
    @(8, 15, 0.152),
    @(8, 16, 0.39),
    @(8, 17, 0.545),
    @(8, 18, 0.231),
    @(8, 19, 0.481),
    @(8, 20, 0.477),
    @(8, 22, 0.458),
    @(8, 23, 0.67),
    @(8, 24, 0.181),
    @(8, 25, 0.426),
    @(8, 26, 0.766),
    @(8, 27, 0.14),
    @(8, 28, 0.374),
    @(8, 29, 0.623),
    @(8, 30, 0.185),
    @(8, 31, 0.43),
    @(8, 32, 0.87),
    @(8, 33, 0.063),
    @(8, 34, 0.25),
    @(8, 35, 0.753),
    @(8, 36, 0.18),
    @(8, 37, 0.425),
    @(8, 38, 0.87),
    @(8, 39, 0.091),
    @(8, 40, 0.302),
    @(8, 41, 0.831),
    @(9, 2, 0.6899999999999999),
    @(9, 3, 0.214),
    @(9, 4, 0.462),
    @(9, 5, 0.476),
    @(9, 8, 0.619),
    @(9, 9, 0.236),
    @(9, 10, 0.486),
    @(9, 11, 0.452),
    @(9, 12, 0.248),
    @(9, 13, 0.498),
    @(9, 14, 0.619),
    @(9, 15, 0.236),
    @(9, 16, 0.486),
    @(9, 17, 0.5),
    @(9, 20, 0.381),
    @(9, 21, 0.236),
    @(9, 22, 0.486),
    @(9, 23, 0.571),
    @(9, 24, 0.245),
    @(9, 25, 0.495),
    @(9, 26, 0.667),
    @(9, 27, 0.222),
    @(9, 28, 0.471),
    @(9, 29, 0.524),
    @(9, 30, 0.249),
    @(9, 31, 0.499),
    @(9, 32, 0.738),
    @(9, 33, 0.193),
    @(9, 34, 0.44),
    @(9, 35, 0.738),
    @(9, 36, 0.193),
    @(9, 37, 0.44),
    @(9, 38, 0.8100000000000001),
    @(9, 39, 0.154),
    @(9, 40, 0.393),
    @(9, 41, 0.762),
    @(10, 2, 0.786),
    @(10, 3, 0.168),
    @(10, 4, 0.41),
    @(10, 5, 0.643),
    @(10, 6, 0.23),
    @(10, 7, 0.479),
    @(10, 8, 0.786),
    @(10, 9, 0.168),
    @(10, 10, 0.41),
    @(10, 11, 0.619),
    @(10, 12, 0.236),
    @(10, 13, 0.486),
    @(10, 14, 0.786),
    @(10, 15, 0.168),
    @(10, 16, 0.41),
    @(10, 17, 0.571),
    @(10, 18, 0.245),
    @(10, 19, 0.495),
    @(10, 20, 0.548),
    @(10, 21, 0.248),
    @(10, 22, 0.498),
    @(10, 23, 0.738),
    @(10, 24, 0.193),
    @(10, 25, 0.44),
    @(10, 26, 0.833),
    @(10, 27, 0.139),
    @(10, 28, 0.373),
    @(10, 29, 0.619),
    @(10, 30, 0.236),
    @(10, 31, 0.486),
    @(10, 32, 0.952),
    @(10, 33, 0.045),
    @(10, 34, 0.213),
    @(10, 35, 0.762),
    @(10, 36, 0.181),
    @(10, 37, 0.426),
    @(10, 38, 0.905),
    @(10, 39, 0.08599999999999999),
    @(10, 40, 0.294),
    @(10, 41, 0.873),
    @(11, 2, 0.833),
    @(11, 3, 0.139),
    @(11, 4, 0.373),
    @(11, 5, 0.714),
    @(11, 6, 0.204),
    @(11, 7, 0.452),
    @(11, 8, 0.857),
    @(11, 9, 0.122),
    @(11, 10, 0.35),
    @(11, 11, 0.619),
    @(11, 12, 0.236),
    @(11, 13, 0.486),
    @(11, 14, 0.8100000000000001),
    @(11, 15, 0.154),
    @(11, 16, 0.393),
    @(11, 17, 0.571),
    @(11, 18, 0.245),
    @(11, 19, 0.495),
    @(11, 20, 0.548),
    @(11, 21, 0.248),
    @(11, 22, 0.498),
    @(11, 23, 0.738),
    @(11, 24, 0.193),
    @(11, 25, 0.44),
    @(11, 26, 0.833),
    @(11, 27, 0.139),
    @(11, 28, 0.373),
    @(11, 29, 0.667),
    @(11, 30, 0.222),
    @(11, 31, 0.471),
    @(11, 32, 0.952),
    @(11, 33, 0.045),
    @(11, 34, 0.213),
    @(11, 35, 0.762),
    @(11, 36, 0.181),
    @(11, 37, 0.426),
    @(11, 38, 0.905),
    @(11, 39, 0.08599999999999999),
    @(11, 40, 0.294),
    @(11, 41, 0.873),
    @(12, 2, 1.314),
    @(12, 3, 0.616),
    @(12, 4, 0.785),
    @(12, 5, 1.633),
    @(12, 6, 1.032),
    @(12, 7, 1.016),
    @(12, 8, 1.556),
    @(12, 9, 1.191),
    @(12, 10, 1.091),
    @(12, 11, 1.423),
    @(12, 12, 0.552),
    @(12, 13, 0.743),
    @(12, 14, 1.353),
    @(12, 15, 0.522),
    @(12, 16, 0.723),
    @(12, 26, 1.257),
    @(12, 27, 0.305),
    @(12, 28, 0.553),
    @(12, 29, 2),
    @(12, 30, 3.935),
    @(12, 31, 1.984),
    @(12, 32, 1.25),
    @(12, 33, 0.237),
    @(12, 34, 0.487),
    @(12, 35, 1.031),
    @(12, 36, 0.03),
    @(12, 37, 0.174),
    @(12, 38, 1.105),
    @(12, 39, 0.094),
    @(12, 40, 0.307),
    @(12, 41, 1.129),
    @(13, 2, 3.429),
    @(13, 3, 1.34),
    @(13, 4, 1.158),
    @(13, 5, 4.553),
    @(13, 6, 0.721),
    @(13, 7, 0.849),
    @(13, 8, 4.524),
    @(13, 9, 0.916),
    @(13, 10, 0.957),
    @(13, 11, 2.297),
    @(13, 12, 0.587),
    @(13, 13, 0.766),
    @(13, 14, 3.333),
    @(13, 15, 0.794),
    @(13, 16, 0.891),
    @(13, 26, 2.725),
    @(13, 27, 3.649),
    @(13, 28, 1.91),
    @(13, 29, 6.244),
    @(13, 30, 2.965),
    @(13, 31, 1.722),
    @(13, 32, 1.643),
    @(13, 33, 0.706),
    @(13, 34, 0.84),
    @(13, 35, 1.31),
    @(13, 36, 0.357),
    @(13, 37, 0.597),
    @(13, 38, 1.714),
    @(13, 39, 0.823),
    @(13, 40, 0.907),
    @(13, 41, 1.556)
)

foreach ($u in $updates) {
    $ws.Cells.Item($u[0], $u[1]).Value = $u[2]
}
